# Add two new worksheets to the opencloning_linkml workbook:
#   - InVivoAssemblySource, inserted right after OverlapExtensionPCRLigationSource
#   - CreLoxRecombinationSource, inserted right after GatewaySource
# Both share the same header row layout as the other simple assembly-source
# sheets (e.g. OverlapExtensionPCRLigationSource):
#   circular | assembly | input | output | type | output_name | database_id | id

$wb = $excel.ActiveWorkbook

$headers = @("circular", "assembly", "input", "output", "type", "output_name", "database_id", "id")

# --- InVivoAssemblySource: goes right after OverlapExtensionPCRLigationSource ---
$afterSheet1 = $wb.Worksheets.Item("OverlapExtensionPCRLigationSource")
$invivo = $wb.Worksheets.Add($null, $afterSheet1)
$invivo.Name = "InVivoAssemblySource"
for ($i = 0; $i -lt $headers.Length; $i++) {
    $invivo.Cells.Item(1, $i + 1).Value = $headers[$i]
}

# --- CreLoxRecombinationSource: goes right after GatewaySource ---
$afterSheet2 = $wb.Worksheets.Item("GatewaySource")
$crelox = $wb.Worksheets.Add($null, $afterSheet2)
$crelox.Name = "CreLoxRecombinationSource"
for ($i = 0; $i -lt $headers.Length; $i++) {
    $crelox.Cells.Item(1, $i + 1).Value = $headers[$i]
}
